# Updated symbol list on Wed Dec 14 04:57:51 UTC 2022 with GitHub Actions
#
# Refreshes the "Price" (column D) values for the coin rows, plus two
# "Worst/Best in 24h" badge updates in column E, to reflect the latest
# scrape of the coinranking data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D (Price) cells are stored as text in the sheet, and several new
# values carry significant trailing zeros (e.g. "3.850", "0.005680") or
# would otherwise render in scientific notation (e.g. "0.00006262") if
# Excel auto-converted them to numbers. Prefixing with an apostrophe forces
# text entry; resetting the style to "Normal" afterwards clears the
# resulting quote-prefix flag so the cell's style matches the original
# (unstyled) cells.
function Set-TextValue($addr, $text) {
    $ws.Range($addr).Value = "'" + $text
    $ws.Range($addr).Style = "Normal"
}

Set-TextValue "D2"  "275.29"
Set-TextValue "D3"  "23.22"
Set-TextValue "D4"  "6.483"
Set-TextValue "D5"  "0.06267"
Set-TextValue "D6"  "3.652"
Set-TextValue "D7"  "6.671"
Set-TextValue "D8"  "1.399"
Set-TextValue "D9"  "0.8314"
Set-TextValue "D10" "0.01381"
Set-TextValue "D11" "0.1625"
Set-TextValue "D12" "0.08288"
Set-TextValue "D13" "0.03436"
Set-TextValue "D14" "0.03121"
Set-TextValue "D15" "0.09304"
Set-TextValue "D16" "3.850"
Set-TextValue "D17" "0.001638"
Set-TextValue "D18" "0.04793"
Set-TextValue "D19" "0.006305"

Set-TextValue "D20" "0.005680"
$ws.Range("E20").Value = "19HotbitTokenHTB"

Set-TextValue "D22" "0.0001500"
Set-TextValue "D26" "0.1249"
Set-TextValue "D41" "0.007045"
Set-TextValue "D42" "0.1163"

Set-TextValue "D43" "0.003350"
$ws.Range("E43").Value = "42CEJICEJIWorstin24h"

Set-TextValue "D44" "0.01215"
Set-TextValue "D45" "0.00006262"
Set-TextValue "D48" "0.7968"
Set-TextValue "D49" "0.01223"
Set-TextValue "D51" "0.01240"

Write-Host "Applied cryptos price refresh"
